$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'ВердиоГаст® Растительный комплекс для улучшения пищеварения (БАД ),  капсулы'
$ws.Range("B2").Value = 80512
$ws.Range("A3").Value = 'Сб. Фитопектол №2 (Грудной сбор №2) 35г'
$ws.Range("B3").Value = 2024
$ws.Range("A4").Value = 'Алтей корни 75г'
$ws.Range("B4").Value = 2128
$ws.Range("A5").Value = 'Аир корневища 75г'
$ws.Range("B5").Value = 3577
$ws.Range("A6").Value = 'Сб. Фитопектол №1 (Грудной сбор №1) 35г'
$ws.Range("B6").Value = 2043
$ws.Range("A7").Value = 'Чага (березовый гриб) 50г'
$ws.Range("B7").Value = 11630
$ws.Range("A8").Value = 'Подорожник большой листья 50г'
$ws.Range("B8").Value = 4088
$ws.Range("A9").Value = 'Липа цветки 35г'
$ws.Range("B9").Value = 10296
$ws.Range("A10").Value = 'Сб. Грудной №4 50г'
$ws.Range("B10").Value = 17808
$ws.Range("A11").Value = 'Чистотел трава 50г'
$ws.Range("B11").Value = 8274
$ws.Range("A12").Value = 'Крапива листья 50г'
$ws.Range("B12").Value = 6935
$ws.Range("A13").Value = 'Багульник болотный побеги 50г'
$ws.Range("B13").Value = 8456
$ws.Range("A14").Value = 'Бессмертник песчаный цветки 30г'
$ws.Range("B14").Value = 16282
$ws.Range("A15").Value = 'Укроп пахучий плоды 50г'
$ws.Range("B15").Value = 40812
$ws.Range("A16").Value = 'Шалфей листья 50г'
$ws.Range("B16").Value = 23800
$ws.Range("A17").Value = 'Солодка корни 50г'
$ws.Range("B17").Value = 21197
$ws.Range("A18").Value = 'Ромашка цветки вн 50г'
$ws.Range("B18").Value = 72894
$ws.Range("A19").Value = 'Дуба кора 75г'
$ws.Range("B19").Value = 49535
$ws.Range("A20").Value = 'Сб. Фитонефрол (Урологический сбор) 50г'
$ws.Range("B20").Value = 12598
$ws.Range("A21").Value = 'Ноготки цветки 50г'
$ws.Range("B21").Value = 17242
$ws.Range("A22").Value = 'Сб. Фитогепатол №2 (Желчегонный сбор №2) 35г'
$ws.Range("B22").Value = 4283
$ws.Range("A23").Value = 'Мята перечная листья 50г'
$ws.Range("B23").Value = 17762
$ws.Range("A24").Value = 'Девясил корневища и корни 50г'
$ws.Range("B24").Value = 13788
$ws.Range("A25").Value = 'Зверобой трава 50г'
$ws.Range("B25").Value = 24584
$ws.Range("A26").Value = 'Лен семена 100г'
$ws.Range("B26").Value = 47613
$ws.Range("A27").Value = 'Мать-и-мачеха листья 35г'
$ws.Range("B27").Value = 22026
$ws.Range("A28").Value = 'Эрва шерстистая трава 30г'
$ws.Range("B28").Value = 15145
$ws.Range("A29").Value = 'Тысячелистник трава 50г'
$ws.Range("B29").Value = 12787
$ws.Range("A30").Value = 'Крушина кора 50г'
$ws.Range("B30").Value = 9582
$ws.Range("A31").Value = 'Полынь горькая трава 50г'
$ws.Range("B31").Value = 37896
$ws.Range("A32").Value = 'Брусника листья 50г'
$ws.Range("B32").Value = 16080
$ws.Range("A33").Value = 'Можжевельник плоды 50г'
$ws.Range("B33").Value = 13886
$ws.Range("A34").Value = 'Пижма цветки 75г'
$ws.Range("B34").Value = 17975
$ws.Range("A35").Value = 'Сенна листья 50г'
$ws.Range("B35").Value = 24751
$ws.Range("A36").Value = 'Боярышник плоды 75г'
$ws.Range("B36").Value = 26426
$ws.Range("A37").Value = 'Эвкалипт прутовидный листья 75г'
$ws.Range("B37").Value = 29002
$ws.Range("A38").Value = 'Береза почки 50г'
$ws.Range("B38").Value = 18141
$ws.Range("A39").Value = 'Чабрец трава 50г'
$ws.Range("B39").Value = 23085
$ws.Range("A40").Value = 'Толокнянка листья 50г'
$ws.Range("B40").Value = 9954
$ws.Range("A41").Value = 'Шиповник плоды низковитаминные 50г'
$ws.Range("B41").Value = 36313
$ws.Range("A42").Value = 'Кукуруза столбики с рыльцами 40г'
$ws.Range("B42").Value = 31465
$ws.Range("A43").Value = 'Валериана корневища с корнями 50г'
$ws.Range("B43").Value = 25108
$ws.Range("A44").Value = 'Спорыш трава 50г'
$ws.Range("B44").Value = 19983
$ws.Range("A45").Value = 'Ламинарии слоевища (морская капуста) 100г'
$ws.Range("B45").Value = 17355
$ws.Range("A46").Value = 'Пустырник трава 50г'
$ws.Range("B46").Value = 18688
$ws.Range("A47").Value = 'Череда трава 50г'
$ws.Range("B47").Value = 23282
$ws.Range("A48").Value = 'Рябина плоды 50г'
$ws.Range("B48").Value = 5362
$ws.Range("A49").Value = 'Фп Детский травяной чай "ФармаЦветик® для иммунитета" 20х1,5 г'
$ws.Range("B49").Value = 1060
$ws.Range("A50").Value = 'Фп Фиточай "Лактафитол" (БАД) 20х1,5 г'
$ws.Range("B50").Value = 10422
$ws.Range("A51").Value = 'Фп Детский травяной чай "ФармаЦветик®  при простуде" 20х1,5 г'
$ws.Range("B51").Value = 2430
$ws.Range("A52").Value = 'Фп Детский травяной чай "ФармаЦветик® для спокойного сна" 20х1,5 г'
$ws.Range("B52").Value = 4600
$ws.Range("A53").Value = 'Фп Детский травяной чай "ФармаЦветик® для животика" 20х1,5 г'
$ws.Range("B53").Value = 4590
$ws.Range("A54").Value = 'Фп "ВердиоГаст® Фиточай для улучшения пищеварения с зеленым чаем"(БАД) 20*1,5г'
$ws.Range("B54").Value = 6650
$ws.Range("A55").Value = 'Фп "ВердиоГаст® Фиточай для улучшения пищеварения с черным чаем" (БАД) 20*1,5г'
$ws.Range("B55").Value = 8680
$ws.Range("A56").Value = 'Фп "Щедрость природы® Фиточай диабетический" 20х2,0 г'
$ws.Range("B56").Value = 612
$ws.Range("A57").Value = 'Фп Сб. Грудной №4 20x2,0г'
$ws.Range("B57").Value = 168296
$ws.Range("A58").Value = 'Фп "Щедрость природы® Фиточай для иммунитета" 20х2,0 г'
$ws.Range("B58").Value = 1566
$ws.Range("A59").Value = 'Фп "Щедрость природы® Фиточай кардиологический" 20х2,0 г'
$ws.Range("B59").Value = 1566
$ws.Range("A60").Value = 'Фп "Щедрость природы® Фиточай при простуде" 20х2,0 г'
$ws.Range("B60").Value = 1332
$ws.Range("A61").Value = 'Фп Сб. Проктофитол (Противогеморроидальный сбор) 20х2,0г'
$ws.Range("B61").Value = 7342
$ws.Range("A62").Value = 'Фп Сб. Фитогастрол (Желудочно-кишечный сбор) 20x2,0г'
$ws.Range("B62").Value = 27909
$ws.Range("A63").Value = 'Фп "Щедрость природы® Фиточай успокоительный"20х2,0 г'
$ws.Range("B63").Value = 2250
$ws.Range("A64").Value = 'Фп Сб. Фитогепатол №3 (Желчегонный сбор №3) 20x2,0г'
$ws.Range("B64").Value = 33137
$ws.Range("A65").Value = 'Фп Сенна листья 20x1,5г'
$ws.Range("B65").Value = 27638
$ws.Range("A66").Value = 'Фп Шиповник плоды 20х2,0г'
$ws.Range("B66").Value = 21546
$ws.Range("A67").Value = 'Фп Чабрец трава 20x1,5 г'
$ws.Range("B67").Value = 30348
$ws.Range("A68").Value = 'Фп Аир корневища 20x1,5г'
$ws.Range("B68").Value = 2183
$ws.Range("A69").Value = 'Фп Брусника листья 20х1,5г'
$ws.Range("B69").Value = 36864
$ws.Range("A70").Value = 'Фп Липа цветки 20x1,5г'
$ws.Range("B70").Value = 35277
$ws.Range("A71").Value = 'Фп Череда трава 20х1,5г'
$ws.Range("B71").Value = 25019
$ws.Range("A72").Value = 'Фп Душица трава 20x1,5г'
$ws.Range("B72").Value = 14202
$ws.Range("A73").Value = 'Фп Золототысячник трава 20х1,5г'
$ws.Range("B73").Value = 2433
$ws.Range("A74").Value = 'Фп Ромашка цветки 20x1,5г'
$ws.Range("B74").Value = 727337
$ws.Range("A75").Value = 'Фп Шалфей листья 20х1,5г'
$ws.Range("B75").Value = 90780
$ws.Range("A76").Value = 'Фп Хвощ полевой трава 20х1,5г'
$ws.Range("B76").Value = 14596
$ws.Range("A77").Value = 'Фп Чистотел трава 20х1,5г'
$ws.Range("B77").Value = 17934
$ws.Range("A78").Value = 'Фп Сб. Желудочный №3 20x2,0г'
$ws.Range("B78").Value = 12150
$ws.Range("A79").Value = 'Фп Сб. Арфазетин-Э 20x2,0г'
$ws.Range("B79").Value = 28691
$ws.Range("A80").Value = 'Фп Мелисса лекарственная трава 20x1,5г'
$ws.Range("B80").Value = 20700
$ws.Range("A81").Value = 'Фп Сб. Фитоседан №3 (Успокоительный сбор №3) 20х2,0г'
$ws.Range("B81").Value = 48903
$ws.Range("A82").Value = 'Фп Зверобой трава 20x1,5г'
$ws.Range("B82").Value = 28949
$ws.Range("A83").Value = 'Фп Сб. Фитонефрол (Урологический сбор) 20x2,0г'
$ws.Range("B83").Value = 127543
$ws.Range("A84").Value = 'Фп Сб. Элекасол 20x2,0г'
$ws.Range("B84").Value = 23454
$ws.Range("A85").Value = 'Фп Пустырник трава 20x1,5г'
$ws.Range("B85").Value = 24330
$ws.Range("A86").Value = 'Фп Сб. Бруснивер 20x2,0г'
$ws.Range("B86").Value = 145435
$ws.Range("A87").Value = 'Фп Сб. Фитоседан №2 (Успокоительный сбор №2) 20x2,0г'
$ws.Range("B87").Value = 33391
$ws.Range("A88").Value = 'Фп Толокнянка листья 20x1,5г'
$ws.Range("B88").Value = 24208
$ws.Range("A89").Value = 'Фп Пастушья сумка трава 20х1,5г'
$ws.Range("B89").Value = 4370
$ws.Range("A90").Value = 'Фп Мята перечная листья 20x1,5г'
$ws.Range("B90").Value = 42161
$ws.Range("A91").Value = 'Фп Крапива листья 20x1,5г'
$ws.Range("B91").Value = 40761
$ws.Range("A92").Value = 'Фп Фиалка трехцветная трава 20x1,5г'
$ws.Range("B92").Value = 2808
$ws.Range("A93").Value = 'Фп Береза листья 20x1,5г'
$ws.Range("B93").Value = 3006
$ws.Range("A94").Value = 'Фп "Щедрость природы® Фиточай очищающий" 20х2,0 г'
$ws.Range("B94").Value = 2556
$ws.Range("A95").Value = 'Фп "Щедрость природы® Фиточай для пищеварения" 20х2,0 г'
$ws.Range("B95").Value = 1512
$ws.Range("A96").Value = 'Фп Пижма цветки 20х1,5г'
$ws.Range("B96").Value = 7230
$ws.Range("A97").Value = 'Фп Подорожник листья 20x1,5г'
$ws.Range("B97").Value = 26552
$ws.Range("A98").Value = 'Фп Ольха соплодия 20х1,5г'
$ws.Range("B98").Value = 3578
$ws.Range("A99").Value = 'Фп Боярышник плоды 20х3,0г'
$ws.Range("B99").Value = 16442
$ws.Range("A100").Value = 'Фп Валериана корневища с корнями 20x1,5г'
$ws.Range("B100").Value = 11852
$ws.Range("A101").Value = 'Фп Ноготки цветки 20x1,5г'
$ws.Range("B101").Value = 25313
$ws.Range("A102").Value = 'Фп Крушина кора 20x1,5г'
$ws.Range("B102").Value = 7289
$ws.Range("A103").Value = 'Фп Тысячелистник трава 20x1,5г'
$ws.Range("B103").Value = 15728
$ws.Range("A104").Value = 'Фп Кровохлебка корневища и корни 20x1,5г'
$ws.Range("B104").Value = 5164
$ws.Range("A105").Value = 'Фп Дуб кора 20х1,5г'
$ws.Range("B105").Value = 7245
$ws.Range("A106").Value = 'Фп Бадан корневища 20x1,5г'
$ws.Range("B106").Value = 2581
$ws.Range("A107").Value = 'Фп Почечный чай листья 20x1,5г'
$ws.Range("B107").Value = 98172
$ws.Range("A108").Value = 'Фп Девясил корневища и корни 20х1,5г'
$ws.Range("B108").Value = 19868
$ws.Range("A109").Value = 'Фп Лапчатка корневища 20x2,5г'
$ws.Range("B109").Value = 9121

$ws.Range("B56").NumberFormat = "0"

$ws.Range("A89").Select()